$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last row (row 61) to shrink the table from 61 to 60 data rows
$ws.Rows.Item(61).Delete()

# Rewrite rows 2-60 with the refreshed dataset values
$ws.Range("A2").Value = 34
$ws.Range("B2").Value = "OPERATIONAL"
$ws.Range("C2").Value = "10 East Arts HUB"
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 8
$ws.Range("A3").Value = 52
$ws.Range("B3").Value = "OPERATIONAL"
$ws.Range("C3").Value = "Academy of Athletic Arts"
$ws.Range("D3").Value = 4.6
$ws.Range("E3").Value = 10
$ws.Range("A4").Value = 15
$ws.Range("B4").Value = "OPERATIONAL"
$ws.Range("C4").Value = "Acceleration Art and Photography"
$ws.Range("D4").Value = 4.6
$ws.Range("E4").Value = 16
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "OPERATIONAL"
$ws.Range("C5").Value = "Art On Main Gallery And Gifts"
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 11
$ws.Range("A6").Value = 55
$ws.Range("B6").Value = "OPERATIONAL"
$ws.Range("C6").Value = "Art To Remember"
$ws.Range("D6").Value = 3.9
$ws.Range("E6").Value = 7
$ws.Range("A7").Value = 12
$ws.Range("B7").Value = "OPERATIONAL"
$ws.Range("C7").Value = "ArtMix"
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 8
$ws.Range("A8").Value = 0
$ws.Range("B8").Value = "OPERATIONAL"
$ws.Range("C8").Value = "Arts Council of Indianapolis"
$ws.Range("D8").Value = 4.8
$ws.Range("E8").Value = 12
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "OPERATIONAL"
$ws.Range("C9").Value = "Arts for Lawrence"
$ws.Range("D9").Value = 4.9
$ws.Range("E9").Value = 19
$ws.Range("A10").Value = 11
$ws.Range("B10").Value = "OPERATIONAL"
$ws.Range("C10").Value = "Arts for Learning"
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("A11").Value = 57
$ws.Range("B11").Value = "OPERATIONAL"
$ws.Range("C11").Value = "Carmel Art Education Studio"
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 12
$ws.Range("A12").Value = 27
$ws.Range("B12").Value = "OPERATIONAL"
$ws.Range("C12").Value = "Carmel International Arts Festival"
$ws.Range("D12").Value = 4.8
$ws.Range("E12").Value = 18
$ws.Range("A13").Value = 22
$ws.Range("B13").Value = "OPERATIONAL"
$ws.Range("C13").Value = "Cat Head Press: Printshop and Artist Cooperative"
$ws.Range("D13").Value = 5
$ws.Range("E13").Value = 26
$ws.Range("A14").Value = 33
$ws.Range("B14").Value = "OPERATIONAL"
$ws.Range("C14").Value = "Clowes Memorial Hall"
$ws.Range("D14").Value = 4.7
$ws.Range("E14").Value = 797
$ws.Range("A15").Value = 5
$ws.Range("B15").Value = "OPERATIONAL"
$ws.Range("C15").Value = "Creative Art Center"
$ws.Range("D15").Value = 4.6
$ws.Range("E15").Value = 5
$ws.Range("A16").Value = 50
$ws.Range("B16").Value = "OPERATIONAL"
$ws.Range("C16").Value = "Creative Energy Arts"
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("A17").Value = 14
$ws.Range("B17").Value = "OPERATIONAL"
$ws.Range("C17").Value = "Creative Expressions Arts"
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("A18").Value = 45
$ws.Range("B18").Value = "OPERATIONAL"
$ws.Range("C18").Value = "Cultural Arts Gallery"
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 1
$ws.Range("A19").Value = 29
$ws.Range("B19").Value = "OPERATIONAL"
$ws.Range("C19").Value = "Dance Arts"
$ws.Range("D19").Value = 4.8
$ws.Range("E19").Value = 10
$ws.Range("A20").Value = 56
$ws.Range("B20").Value = "OPERATIONAL"
$ws.Range("C20").Value = "Dance Magic Performing Arts"
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("A21").Value = 32
$ws.Range("B21").Value = "OPERATIONAL"
$ws.Range("C21").Value = "Edison School Of The Arts"
$ws.Range("D21").Value = 3.6
$ws.Range("E21").Value = 36
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "OPERATIONAL"
$ws.Range("C22").Value = "Gallery 924"
$ws.Range("D22").Value = 4.9
$ws.Range("E22").Value = 15
$ws.Range("A23").Value = 54
$ws.Range("B23").Value = "OPERATIONAL"
$ws.Range("C23").Value = "Glass Arts Indiana, Inc."
$ws.Range("D23").Value = 4.8
$ws.Range("E23").Value = 16
$ws.Range("A24").Value = 44
$ws.Range("B24").Value = "OPERATIONAL"
$ws.Range("C24").Value = "H.J. Ricks Centre For The Arts"
$ws.Range("D24").Value = 4.5
$ws.Range("E24").Value = 44
$ws.Range("A25").Value = 16
$ws.Range("B25").Value = "OPERATIONAL"
$ws.Range("C25").Value = "Hancock County Arts Council, Twenty North Gallery"
$ws.Range("D25").Value = 4.7
$ws.Range("E25").Value = 3
$ws.Range("A26").Value = 18
$ws.Range("B26").Value = "OPERATIONAL"
$ws.Range("C26").Value = "Harrison Center"
$ws.Range("D26").Value = 4.8
$ws.Range("E26").Value = 134
$ws.Range("A27").Value = 37
$ws.Range("B27").Value = "OPERATIONAL"
$ws.Range("C27").Value = "Herron School of Art and Design"
$ws.Range("D27").Value = 4.8
$ws.Range("E27").Value = 20
$ws.Range("A28").Value = 35
$ws.Range("B28").Value = "OPERATIONAL"
$ws.Range("C28").Value = "Herron School of Art and Design - Eskenazi Fine Arts Center"
$ws.Range("D28").Value = 4.5
$ws.Range("E28").Value = 2
$ws.Range("A29").Value = 13
$ws.Range("B29").Value = "OPERATIONAL"
$ws.Range("C29").Value = "High Frequency Arts"
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = 3
$ws.Range("A30").Value = 46
$ws.Range("B30").Value = "OPERATIONAL"
$ws.Range("C30").Value = "Ignition Arts, LLC"
$ws.Range("D30").Value = 4.9
$ws.Range("E30").Value = 8
$ws.Range("A31").Value = 2
$ws.Range("B31").Value = "OPERATIONAL"
$ws.Range("C31").Value = "Indiana Arts Commission"
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = "OPERATIONAL"
$ws.Range("C32").Value = "Indiana Fine Arts Academy"
$ws.Range("D32").Value = 4
$ws.Range("E32").Value = 4
$ws.Range("A33").Value = 41
$ws.Range("B33").Value = "OPERATIONAL"
$ws.Range("C33").Value = "Indiana Performing Arts Centre"
$ws.Range("D33").Value = 4.8
$ws.Range("E33").Value = 6
$ws.Range("A34").Value = 3
$ws.Range("B34").Value = "OPERATIONAL"
$ws.Range("C34").Value = "Indianapolis Art Center"
$ws.Range("D34").Value = 4.7
$ws.Range("E34").Value = 112
$ws.Range("A35").Value = 59
$ws.Range("B35").Value = "OPERATIONAL"
$ws.Range("C35").Value = "Indianapolis Arts Chorale"
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("A36").Value = 19
$ws.Range("B36").Value = "OPERATIONAL"
$ws.Range("C36").Value = "Indianapolis Artsgarden"
$ws.Range("D36").Value = 4.6
$ws.Range("E36").Value = 146
$ws.Range("A37").Value = 51
$ws.Range("B37").Value = "OPERATIONAL"
$ws.Range("C37").Value = "Jazz Arts Society of In Inc"
$ws.Range("D37").Value = 4
$ws.Range("E37").Value = 1
$ws.Range("A38").Value = 21
$ws.Range("B38").Value = "OPERATIONAL"
$ws.Range("C38").Value = "Korka International Arts"
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0
$ws.Range("A39").Value = 42
$ws.Range("B39").Value = "OPERATIONAL"
$ws.Range("C39").Value = "Landmark Keystone Art Cinema"
$ws.Range("D39").Value = 4.4
$ws.Range("E39").Value = 513
$ws.Range("A40").Value = 10
$ws.Range("B40").Value = "OPERATIONAL"
$ws.Range("C40").Value = "Magdalena Gallery of Arts"
$ws.Range("D40").Value = 5
$ws.Range("E40").Value = 1
$ws.Range("A41").Value = 48
$ws.Range("B41").Value = "OPERATIONAL"
$ws.Range("C41").Value = "Midland Arts & Antiques Market"
$ws.Range("D41").Value = 4.6
$ws.Range("E41").Value = 688
$ws.Range("A42").Value = 24
$ws.Range("B42").Value = "OPERATIONAL"
$ws.Range("C42").Value = "Newfields"
$ws.Range("D42").Value = 4.7
$ws.Range("E42").Value = 3896
$ws.Range("A43").Value = 9
$ws.Range("B43").Value = "OPERATIONAL"
$ws.Range("C43").Value = "Nickel Plate Arts"
$ws.Range("D43").Value = 4.9
$ws.Range("E43").Value = 41
$ws.Range("A44").Value = 47
$ws.Range("B44").Value = "OPERATIONAL"
$ws.Range("C44").Value = "Performing Arts"
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("A45").Value = 39
$ws.Range("B45").Value = "OPERATIONAL"
$ws.Range("C45").Value = "Phiri Art"
$ws.Range("D45").Value = 4.5
$ws.Range("E45").Value = 4
$ws.Range("A46").Value = 58
$ws.Range("B46").Value = "OPERATIONAL"
$ws.Range("C46").Value = "Red Barn Arts Collective"
$ws.Range("D46").Value = 5
$ws.Range("E46").Value = 1
$ws.Range("A47").Value = 17
$ws.Range("B47").Value = "OPERATIONAL"
$ws.Range("C47").Value = "Schrott Center for the Arts"
$ws.Range("D47").Value = 4.8
$ws.Range("E47").Value = 106
$ws.Range("A48").Value = 1
$ws.Range("B48").Value = "OPERATIONAL"
$ws.Range("C48").Value = "Sho Arts"
$ws.Range("D48").Value = 5
$ws.Range("E48").Value = 1
$ws.Range("A49").Value = 31
$ws.Range("B49").Value = "OPERATIONAL"
$ws.Range("C49").Value = "Steve Haigh Fine Art"
$ws.Range("D49").Value = 5
$ws.Range("E49").Value = 1
$ws.Range("A50").Value = 26
$ws.Range("B50").Value = "OPERATIONAL"
$ws.Range("C50").Value = "Studio Alchemy LLC"
$ws.Range("D50").Value = 5
$ws.Range("E50").Value = 5
$ws.Range("A51").Value = 36
$ws.Range("B51").Value = "OPERATIONAL"
$ws.Range("C51").Value = "Ten West Center for the Arts"
$ws.Range("D51").Value = 4.8
$ws.Range("E51").Value = 4
$ws.Range("A52").Value = 25
$ws.Range("B52").Value = "OPERATIONAL"
$ws.Range("C52").Value = "The Art Studio of Carmel"
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 0
$ws.Range("A53").Value = 6
$ws.Range("B53").Value = "OPERATIONAL"
$ws.Range("C53").Value = "The Carmel Arts Council Children’s Art Gallery"
$ws.Range("D53").Value = 5
$ws.Range("E53").Value = 1
$ws.Range("A54").Value = 23
$ws.Range("B54").Value = "OPERATIONAL"
$ws.Range("C54").Value = "The Center for the Performing Arts"
$ws.Range("D54").Value = 4.7
$ws.Range("E54").Value = 1193
$ws.Range("A55").Value = 7
$ws.Range("B55").Value = "OPERATIONAL"
$ws.Range("C55").Value = "The Murphy Art Center"
$ws.Range("D55").Value = 5
$ws.Range("E55").Value = 1
$ws.Range("A56").Value = 53
$ws.Range("B56").Value = "OPERATIONAL"
$ws.Range("C56").Value = "The Palladium at the Center for the Performing Arts"
$ws.Range("D56").Value = 4.9
$ws.Range("E56").Value = 197
$ws.Range("A57").Value = 28
$ws.Range("B57").Value = "OPERATIONAL"
$ws.Range("C57").Value = "United Art & Education"
$ws.Range("D57").Value = 4.5
$ws.Range("E57").Value = 132
$ws.Range("A58").Value = 38
$ws.Range("B58").Value = "OPERATIONAL"
$ws.Range("C58").Value = "Warren Performing Arts Center"
$ws.Range("D58").Value = 4.5
$ws.Range("E58").Value = 374
$ws.Range("A59").Value = 49
$ws.Range("B59").Value = "OPERATIONAL"
$ws.Range("C59").Value = "Zionsville Performing Arts Center"
$ws.Range("D59").Value = 4.7
$ws.Range("E59").Value = 175
$ws.Range("A60").Value = 40
$ws.Range("B60").Value = "OPERATIONAL"
$ws.Range("C60").Value = "carmel academy of the arts"
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
